# Saldo.xlsx edit:
#  - Remove the rows for ROBERIO (004586209/24000), CEZAR (004482090/11779.81)
#    and CLINEO (004204344/1100).
#  - Re-order the block of rows around "HFR" (004361159) so that
#    LAURA/ELIANE/JOSE/THIAGO/THOMAS now come right after DOUGLAS, and the
#    HFR row moves to the end of that block with its balance updated to
#    14873.25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three obsolete rows bottom-up so earlier row numbers stay valid.
$ws.Rows.Item(14).Delete()   # 004204344 CLINEO   1100
$ws.Rows.Item(11).Delete()   # 004482090 CEZAR    11779.81
$ws.Rows.Item(9).Delete()    # 004586209 ROBERIO  24000

# After the deletions above, rows 4-9 hold (in order):
#   4 HFR   5 LAURA   6 ELIANE   7 JOSE   8 THIAGO   9 THOMAS
# Grab that block, then rotate it so HFR lands last (row 9) with its new
# balance, and everything else shifts up by one row.
$block = $ws.Range("A4:C9").Value2

for ($i = 4; $i -le 8; $i++) {
    $srcRow = $i - 4 + 2   # rows 5..9 of the block (LAURA..THOMAS)
    $ws.Cells.Item($i, 1).NumberFormat = "@"
    $ws.Cells.Item($i, 1).Value = $block[$srcRow, 1]
    $ws.Cells.Item($i, 2).Value = $block[$srcRow, 2]
    $ws.Cells.Item($i, 3).Value = $block[$srcRow, 3]
}

$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = $block[1, 1]
$ws.Cells.Item(9, 2).Value = $block[1, 2]
$ws.Cells.Item(9, 3).Value = 14873.25
